$d = $word.ActiveDocument

$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="63BE0E92" w14:textId="2634EE53" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:proofErr w:type="spellStart"/><w:r><w:t>Woohoo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0062343C"><w:t>Docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Templater</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is ||adjective||</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(2).Range.InsertXML($xml2)

$xml4 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1B9CAA6E" w14:textId="77777777" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:r><w:t>||</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>foo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>||</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(4).Range.InsertXML($xml4)

$xml5 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6442F136" w14:textId="77777777" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>||</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>key</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>||</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(5).Range.InsertXML($xml5)

$xml7 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="02609025" w14:textId="77777777" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:pPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>header</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(7).Range.InsertXML($xml7)

$xml9 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3C42B4C6" w14:textId="77777777" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:r><w:t>||</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>foo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>||</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(9).Range.InsertXML($xml9)

$xml10 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="661DDBC1" w14:textId="77777777" w:rsidR="009E187F" w:rsidRDefault="009E187F" w:rsidP="009E187F"><w:pPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>||</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>key</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>||</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(10).Range.InsertXML($xml10)

$xml25 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="025EC860" w14:textId="0CD3EF56" w:rsidR="00BF0AF0" w:rsidRDefault="00BF0AF0" w:rsidP="00217966"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="003A01A4"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>||</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>total</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>_loan_amount_</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>currency_</w:t></w:r><w:r w:rsidR="006D003D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>words</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D003D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">|| </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(25).Range.InsertXML($xml25)

$xml26 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="45E95AFC" w14:textId="0939B7FC" w:rsidR="00217966" w:rsidRDefault="00217966" w:rsidP="00217966"><w:r><w:t>It can even replace placeholders wit</w:t></w:r><w:r w:rsidRPr="009E187F"><w:rPr><w:b/></w:rPr><w:t>h ||</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009E187F"><w:rPr><w:b/></w:rPr><w:t>supe</w:t></w:r><w:r w:rsidRPr="009E187F"><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>r_adject</w:t></w:r><w:r w:rsidRPr="009E187F"><w:rPr><w:b/></w:rPr><w:t>ive</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>|</w:t></w:r><w:r w:rsidRPr="009E187F"><w:rPr><w:i/></w:rPr><w:t>| f</w:t></w:r><w:r><w:t>ormatting</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(26).Range.InsertXML($xml26)

$xml27 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="448F3B79" w14:textId="77777777" w:rsidR="00217966" w:rsidRPr="009E187F" w:rsidRDefault="00217966" w:rsidP="009E187F"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(27).Range.InsertXML($xml27)
